$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the slightly drifted timestamp value in A3 (automatic scheduled task update)
$ws.Range("A3").Value = 45864.08353502315

# Append the new row of scheduled-task data (row 4)
$ws.Range("A4").Value = 45864.3336263014
$ws.Range("A4").NumberFormat = $ws.Range("A3").NumberFormat

$ws.Range("B4").Value = 2025
$ws.Range("C4").Value = 30
$ws.Range("D4").Value = 12.83
$ws.Range("E4").Value = 92
$ws.Range("F4").Value = 80.86
$ws.Range("G4").Value = 10.07
$ws.Range("H4").Value = "WNW"
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = "08:00:25"
